$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F "想去人数" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12809
$ws1.Range("F3").Value = 622
$ws1.Range("F5").Value = 30
$ws1.Range("F6").Value = 317
$ws1.Range("F8").Value = 235
$ws1.Range("F9").Value = 12805
$ws1.Range("F10").Value = 38
$ws1.Range("F11").Value = 19
$ws1.Range("F12").Value = 5215
$ws1.Range("F15").Value = 12
$ws1.Range("F16").Value = 28
$ws1.Range("F18").Value = 34
$ws1.Range("F20").Value = 673
$ws1.Range("F21").Value = 2849
$ws1.Range("F22").Value = 6167
$ws1.Range("F24").Value = 3621
$ws1.Range("F26").Value = 42

# Sheet "演出" (Performances) - column F update
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 21

# Sheet "全部类型" (All Types) - column F "想去人数" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12809
$ws4.Range("F3").Value = 622
$ws4.Range("F5").Value = 30
$ws4.Range("F6").Value = 317
$ws4.Range("F7").Value = 21
$ws4.Range("F9").Value = 235
$ws4.Range("F10").Value = 12805
$ws4.Range("F11").Value = 38
$ws4.Range("F12").Value = 19
$ws4.Range("F13").Value = 5215
$ws4.Range("F16").Value = 12
$ws4.Range("F17").Value = 28
$ws4.Range("F19").Value = 34
$ws4.Range("F21").Value = 673
$ws4.Range("F22").Value = 2849
$ws4.Range("F24").Value = 6167
$ws4.Range("F26").Value = 3621
$ws4.Range("F28").Value = 42

$wb.Save()
